$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'freuen',
    'formen',
    'stecken',
    'wirken',
    'seufzen',
    'stehlen',
    'arten',
    'lesen',
    'bluten',
    'knarren',
    'sperren',
    'bellen',
    'spinnen',
    'folgen',
    'dringen',
    'lockern',
    'runden',
    'treiben',
    'erben',
    'backen',
    'äußern',
    'spielen',
    'jubeln',
    'jagen',
    'machen',
    'decken',
    'feiern',
    'trauen',
    'grüßen',
    'mauern',
    'münzen',
    'schwächen',
    'schrecken',
    'töten',
    'warnen',
    'schmecken',
    'flüchten',
    'rufen',
    'werfen',
    'schulden',
    'kehren',
    'malen',
    'wehtun',
    'sprengen',
    'suchen',
    'kichern',
    'filmen',
    'heulen',
    'klingen',
    'schwören',
    'ärgern',
    'ehren',
    'achten',
    'irren',
    'reizen',
    'sorgen',
    'kümmern',
    'weichen',
    'zögern',
    'kosten',
    'biegen',
    'fischen',
    'werden',
    'schlucken',
    'brauchen',
    'schenken',
    'schreiten',
    'heilen',
    'saufen',
    'bitten',
    'hören',
    'wenden',
    'geben',
    'führen',
    'fließen',
    'zielen',
    'sterben',
    'bauen',
    'wüten',
    'fällen',
    'stammen',
    'heben',
    'tollen',
    'dienen',
    'altern',
    'trennen',
    'fangen',
    'streichen',
    'lügen',
    'schwingen',
    'planen',
    'graben',
    'hauen',
    'scheinen',
    'helfen',
    'betteln',
    'tropfen',
    'fallen',
    'sinken',
    'mögen',
    'scheitern',
    'greifen',
    'wundern',
    'fahren',
    'spüren',
    'boxen',
    'rasen',
    'siegen',
    'ändern',
    'liefern',
    'flehen',
    'pfeifen',
    'räumen',
    'sichern',
    'klettern',
    'bergen',
    'enden',
    'zünden',
    'gründen',
    'liegen',
    'gelten',
    'zeigen',
    'drehen',
    'platzen',
    'quälen',
    'wachsen',
    'loben',
    'pflanzen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $words[$i]
}
